$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates keyed by row number => Price (D) / Volume(1h) (E) text.
# Only columns actually changed per row are included; null entries are skipped.
$updates = @{
    2  = @{ D = "27.517.22";  E = "  -1.47%  " }
    3  = @{ D = "1.842.10";   E = "  -2.05%  " }
    4  = @{ D = "1.006";      E = "  -1.32%  " }
    5  = @{ D = "334.33";     E = "  -0.07%  " }
    6  = @{ E = "  -1.17%  " }
    7  = @{ D = "0.4628";     E = "  -0.76%  " }
    8  = @{ D = "0.3847";     E = "  -1.43%  " }
    9  = @{ D = "45.84";      E = "  -2.94%  " }
    10 = @{ D = "0.07893";    E = "  -0.62%  " }
    11 = @{ D = "0.9969";     E = "  -0.79%  " }
    12 = @{ D = "21.42";      E = "  -0.60%  " }
    13 = @{ E = "  +0.65%  " }
    14 = @{ D = "1.840.39";   E = "  -3.97%  " }
    15 = @{ D = "7.119";      E = "  +0.79%  " }
    16 = @{ D = "1.008";      E = "  -1.29%  " }
    17 = @{ D = "88.29";      E = "  +1.67%  " }
    18 = @{ D = "0.06671";    E = "  -1.25%  " }
    19 = @{ E = "  -0.68%  " }
    20 = @{ D = "17.10";      E = "  +0.75%  " }
    21 = @{ E = "  -1.25%  " }
    22 = @{ D = "27.507.00";  E = "  -1.60%  " }
    23 = @{ D = "5.378";      E = "  -1.52%  " }
    24 = @{ D = "10.87";      E = "  -0.19%  " }
    25 = @{ D = "2.314";      E = "  -1.61%  " }
    26 = @{ D = "158.85";     E = "  -0.61%  " }
    27 = @{ D = "2.058.27";   E = "  -3.74%  " }
    28 = @{ D = "19.47" }
    29 = @{ E = "  +2.33%  " }
    30 = @{ D = "5.394";      E = "  -0.86%  " }
    31 = @{ D = "119.74";     E = "  -0.96%  " }
    32 = @{ D = "0.9742";     E = "  +2.15%  " }
    33 = @{ D = "0.09395";    E = "  -0.83%  " }
    34 = @{ D = "3.594";      E = "  -1.92%  " }
    35 = @{ D = "5.293";      E = "  -0.03%  " }
    36 = @{ E = "  -1.57%  " }
    37 = @{ D = "0.06015";    E = "  -1.51%  " }
    38 = @{ D = "0.02223";    E = "  -0.18%  " }
    39 = @{ D = "8.261";      E = "  +2.01%  " }
    40 = @{ D = "1.177";      E = "  -2.84%  " }
    41 = @{ D = "0.5882";     E = "  +0.04%  " }
    44 = @{ E = "  -2.40%  " }
    45 = @{ D = "0.5574";     E = "  -0.75%  " }
    46 = @{ D = "12.10";      E = "  +0.27%  " }
    47 = @{ E = "  -0.16%  " }
    48 = @{ D = "0.06685";    E = "  -2.84%  " }
    49 = @{ D = "110.14";     E = "  -2.79%  " }
    50 = @{ D = "1.046";      E = "  -1.21%  " }
    51 = @{ D = "1.007";      E = "  -1.24%  " }
}

# Capture the default (un-styled) cell style from an untouched data cell so it
# can be restored after the temporary "@" text format forces these numeric-
# looking strings (e.g. "1.006") to be stored as text instead of numbers.
$defaultStyle = $ws.Range("D4").Style

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals.D
        $cell.Style = $defaultStyle
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals.E
        $cell.Style = $defaultStyle
    }
}

# Rows 42 and 43 swap coin identity (Aptos <-> Algorand) along with new values.
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1857"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E42").Style = $defaultStyle

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.32"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("E43").Style = $defaultStyle
